$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 178, shifting existing rows 178:269 down to 179:270.
$ws.Rows("178:178").Insert()

# Populate the newly inserted row 178 with the new weekly observation.
$ws.Cells.Item(178, 1).Value = 9
$ws.Cells.Item(178, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(178, 3).Value = "Metropolitana"
$ws.Cells.Item(178, 4).Value = "2022-09-09"
$ws.Cells.Item(178, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(178, 5).Value = 13
$ws.Cells.Item(178, 6).Value = 100112026
$ws.Cells.Item(178, 7).Value = "Haba"
$ws.Cells.Item(178, 8).Value = "Sin especificar"
$ws.Cells.Item(178, 9).Value = "Primera"
$ws.Cells.Item(178, 10).Value = 70
$ws.Cells.Item(178, 11).Value = 12000
$ws.Cells.Item(178, 12).Value = 13000
$ws.Cells.Item(178, 13).Value = 12429
$ws.Cells.Item(178, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(178, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(178, 16).Value = 497
$ws.Cells.Item(178, 17).Value = 25
$ws.Cells.Item(178, 18).Value = "Hortaliza"
